# Updates cryptos list prices / 1h volume percentages (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Scratch cell used to push new values through Excel as literal text
# (PasteSpecial values-only), so numeric-looking strings such as "604.23"
# are not silently converted into real numbers and no new cell styles
# (number formats) get introduced along the way.
$helper = $ws.Range("Z1")

$helper.Formula = "=""66.891.59"""
$helper.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$helper.Formula = "=""  +2.91%  """
$helper.Copy()
$ws.Range("E2").PasteSpecial(-4163)

$helper.Formula = "=""3.201.06"""
$helper.Copy()
$ws.Range("D3").PasteSpecial(-4163)

$helper.Formula = "=""  +1.70%  """
$helper.Copy()
$ws.Range("E3").PasteSpecial(-4163)

$helper.Formula = "=""  +0.00%  """
$helper.Copy()
$ws.Range("E4").PasteSpecial(-4163)

$helper.Formula = "=""604.23"""
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$helper.Formula = "=""  +4.25%  """
$helper.Copy()
$ws.Range("E5").PasteSpecial(-4163)

$helper.Formula = "=""157.26"""
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)

$helper.Formula = "=""  +5.64%  """
$helper.Copy()
$ws.Range("E6").PasteSpecial(-4163)

$helper.Formula = "=""  +0.01%  """
$helper.Copy()
$ws.Range("E7").PasteSpecial(-4163)

$helper.Formula = "=""0.557"""
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)

$helper.Formula = "=""  +6.24%  """
$helper.Copy()
$ws.Range("E8").PasteSpecial(-4163)

$helper.Formula = "=""3.199.01"""
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)

$helper.Formula = "=""  +1.66%  """
$helper.Copy()
$ws.Range("E9").PasteSpecial(-4163)

$helper.Formula = "=""  +1.87%  """
$helper.Copy()
$ws.Range("E10").PasteSpecial(-4163)

$helper.Formula = "=""  -3.71%  """
$helper.Copy()
$ws.Range("E11").PasteSpecial(-4163)

$helper.Formula = "=""  +3.50%  """
$helper.Copy()
$ws.Range("E12").PasteSpecial(-4163)

$helper.Formula = "=""0.0000269"""
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)

$helper.Formula = "=""  +1.99%  """
$helper.Copy()
$ws.Range("E13").PasteSpecial(-4163)

$helper.Formula = "=""39.24"""
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)

$helper.Formula = "=""  +5.65%  """
$helper.Copy()
$ws.Range("E14").PasteSpecial(-4163)

$helper.Formula = "=""3.725.01"""
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)

$helper.Formula = "=""  +1.76%  """
$helper.Copy()
$ws.Range("E15").PasteSpecial(-4163)

$helper.Formula = "=""66.768.74"""
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)

$helper.Formula = "=""  +2.86%  """
$helper.Copy()
$ws.Range("E16").PasteSpecial(-4163)

$helper.Formula = "=""7.49"""
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)

$helper.Formula = "=""  +4.98%  """
$helper.Copy()
$ws.Range("E17").PasteSpecial(-4163)

$helper.Formula = "=""3.202.78"""
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)

$helper.Formula = "=""  +1.99%  """
$helper.Copy()
$ws.Range("E18").PasteSpecial(-4163)

$helper.Formula = "=""  +0.91%  """
$helper.Copy()
$ws.Range("E19").PasteSpecial(-4163)

$helper.Formula = "=""522.13"""
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)

$helper.Formula = "=""  +3.60%  """
$helper.Copy()
$ws.Range("E20").PasteSpecial(-4163)

$helper.Formula = "=""15.50"""
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)

$helper.Formula = "=""  +2.88%  """
$helper.Copy()
$ws.Range("E21").PasteSpecial(-4163)

$helper.Formula = "=""0.744"""
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)

$helper.Formula = "=""  +4.16%  """
$helper.Copy()
$ws.Range("E22").PasteSpecial(-4163)

$helper.Formula = "=""8.22"""
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)

$helper.Formula = "=""  +6.38%  """
$helper.Copy()
$ws.Range("E23").PasteSpecial(-4163)

$helper.Formula = "=""15.11"""
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)

$helper.Formula = "=""  -0.29%  """
$helper.Copy()
$ws.Range("E24").PasteSpecial(-4163)

$helper.Formula = "=""85.48"""
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)

$helper.Formula = "=""  +1.50%  """
$helper.Copy()
$ws.Range("E25").PasteSpecial(-4163)

$helper.Formula = "=""  -0.09%  """
$helper.Copy()
$ws.Range("E26").PasteSpecial(-4163)

$helper.Formula = "=""9.28"""
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)

$helper.Formula = "=""  +1.81%  """
$helper.Copy()
$ws.Range("E27").PasteSpecial(-4163)

$helper.Formula = "=""3.02"""
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)

$helper.Formula = "=""  +3.72%  """
$helper.Copy()
$ws.Range("E28").PasteSpecial(-4163)

$helper.Formula = "=""  +10.12%  """
$helper.Copy()
$ws.Range("E29").PasteSpecial(-4163)

$helper.Formula = "=""3.02"""
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)

$helper.Formula = "=""  +8.27%  """
$helper.Copy()
$ws.Range("E30").PasteSpecial(-4163)

$helper.Formula = "=""7.03"""
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)

$helper.Formula = "=""  +9.36%  """
$helper.Copy()
$ws.Range("E31").PasteSpecial(-4163)

$helper.Formula = "=""28.34"""
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)

$helper.Formula = "=""  +3.07%  """
$helper.Copy()
$ws.Range("E32").PasteSpecial(-4163)

$helper.Formula = "=""  +2.98%  """
$helper.Copy()
$ws.Range("E33").PasteSpecial(-4163)

$helper.Formula = "=""  +0.12%  """
$helper.Copy()
$ws.Range("E34").PasteSpecial(-4163)

$helper.Formula = "=""6.58"""
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)

$helper.Formula = "=""  +1.52%  """
$helper.Copy()
$ws.Range("E35").PasteSpecial(-4163)

$helper.Formula = "=""524.18"""
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)

$helper.Formula = "=""  +10.22%  """
$helper.Copy()
$ws.Range("E36").PasteSpecial(-4163)

$helper.Formula = "=""55.22"""
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)

$helper.Formula = "=""  +0.48%  """
$helper.Copy()
$ws.Range("E37").PasteSpecial(-4163)

$helper.Formula = "=""0.0905"""
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)

$helper.Formula = "=""  +2.10%  """
$helper.Copy()
$ws.Range("E38").PasteSpecial(-4163)

$helper.Formula = "=""0.0428"""
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)

$helper.Formula = "=""  +3.54%  """
$helper.Copy()
$ws.Range("E39").PasteSpecial(-4163)

$helper.Formula = "=""0.127"""
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)

$helper.Formula = "=""  +8.87%  """
$helper.Copy()
$ws.Range("E40").PasteSpecial(-4163)

$helper.Formula = "=""8.92"""
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)

$helper.Formula = "=""  +1.97%  """
$helper.Copy()
$ws.Range("E41").PasteSpecial(-4163)

$helper.Formula = "=""  -0.75%  """
$helper.Copy()
$ws.Range("E42").PasteSpecial(-4163)

$helper.Formula = "=""  +15.50%  """
$helper.Copy()
$ws.Range("E43").PasteSpecial(-4163)

$helper.Formula = "=""0.302"""
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)

$helper.Formula = "=""  +7.25%  """
$helper.Copy()
$ws.Range("E44").PasteSpecial(-4163)

$helper.Formula = "=""  +2.52%  """
$helper.Copy()
$ws.Range("E45").PasteSpecial(-4163)

$helper.Formula = "=""2.899.32"""
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)

$helper.Formula = "=""  -3.33%  """
$helper.Copy()
$ws.Range("E46").PasteSpecial(-4163)

$helper.Formula = "=""28.67"""
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)

$helper.Formula = "=""  +1.27%  """
$helper.Copy()
$ws.Range("E47").PasteSpecial(-4163)

$helper.Formula = "=""  +11.09%  """
$helper.Copy()
$ws.Range("E48").PasteSpecial(-4163)

$helper.Formula = "=""  +3.78%  """
$helper.Copy()
$ws.Range("E49").PasteSpecial(-4163)

$helper.Formula = "=""  -0.01%  """
$helper.Copy()
$ws.Range("E50").PasteSpecial(-4163)

$helper.Formula = "=""  +4.38%  """
$helper.Copy()
$ws.Range("E51").PasteSpecial(-4163)

$helper.ClearContents()
$excel.CutCopyMode = $false
